$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The table is being reshaped from a 6-col x 2-row data table (A1:F3) into a
# 3-col x 7-row table (A1:C7): a 3-cell header row followed by six single
# -column data rows.
# ---------------------------------------------------------------------------

# The existing data cells A2, B2, C2, D2, E2, F2 already hold the exact text
# values the new column B needs ("5,206", "2,204", "3,002", "60", "30", "30"),
# so relocate them (copy+paste VALUES, not a re-typed literal) instead of
# retyping, which keeps them stored as text/shared-strings without minting
# any new number-format styles. Walk back-to-front so a source cell is
# always read before anything overwrites it.
$ws.Range("F2").Copy()
$ws.Range("B7").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("E2").Copy()
$ws.Range("B6").PasteSpecial(-4163)
$ws.Range("D2").Copy()
$ws.Range("B5").PasteSpecial(-4163)
$ws.Range("C2").Copy()
$ws.Range("B4").PasteSpecial(-4163)
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4163)

# Drop the now-stale source cells (old row 2/3 leftovers outside column B).
$ws.Range("A2:A3").Clear()
$ws.Range("C2:F2").Clear()

# Header row. A1 already carries the bold/border/centered style (style index
# 1) - keep it and just retitle it, then clone that exact formatting onto the
# two new header cells B1/C1.
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)  # xlPasteFormats

# "2019" must stay textual (not become the number 2019). A formula whose
# result is a literal string evaluates to a text cell without requiring any
# text number-format, so stage it in a scratch cell via a formula, then
# bring only the VALUE over into B1 (keeping the header formatting already
# pasted above instead of any scratch formatting).
$scratch = $ws.Range("Z100")
$scratch.Formula = '="2019"'
$scratch.Copy()
$ws.Range("B1").PasteSpecial(-4163)   # xlPasteValues
$scratch.Clear()

$ws.Range("C1").Value = "Unnamed: 1"

$ws.Range("A1:C1").Font.Bold = $true
